# Scheduled runner update: refresh currentAveragePrice/Profit columns (H-N)
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the
# latest market-board derived figures.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Cells.Item(13, 8).Value = 43333.332
$ws_ALC.Cells.Item(13, 10).Value = 40000
$ws_ALC.Cells.Item(13, 12).Value = 40000
$ws_ALC.Cells.Item(13, 14).Value = -40338
$ws_ALC.Cells.Item(39, 8).Value = 267.55554
$ws_ALC.Cells.Item(39, 9).Value = 99.75
$ws_ALC.Cells.Item(39, 11).Value = 299.25
$ws_ALC.Cells.Item(39, 13).Value = -3.25
$ws_ALC.Cells.Item(61, 8).Value = 475.125
$ws_ALC.Cells.Item(61, 9).Value = 475.125
$ws_ALC.Cells.Item(61, 11).Value = 1425.375
$ws_ALC.Cells.Item(61, 13).Value = -1253.375
$ws_ALC.Cells.Item(99, 8).Value = 771.8
$ws_ALC.Cells.Item(99, 9).Value = 459.66666
$ws_ALC.Cells.Item(99, 10).Value = 1240
$ws_ALC.Cells.Item(99, 11).Value = 1378.99998
$ws_ALC.Cells.Item(99, 12).Value = 3720
$ws_ALC.Cells.Item(99, 13).Value = 119.0000199999999
$ws_ALC.Cells.Item(99, 14).Value = -6716
$ws_ALC.Cells.Item(113, 8).Value = 2087.6
$ws_ALC.Cells.Item(113, 9).Value = 2087.6
$ws_ALC.Cells.Item(113, 11).Value = 2087.6
$ws_ALC.Cells.Item(113, 13).Value = 1166.4
$ws_ALC.Cells.Item(133, 8).Value = 25703.75
$ws_ALC.Cells.Item(133, 10).Value = 25703.75
$ws_ALC.Cells.Item(133, 12).Value = 25703.75
$ws_ALC.Cells.Item(133, 14).Value = -35823.75
$ws_ALC.Cells.Item(138, 8).Value = 6533515.5
$ws_ALC.Cells.Item(138, 9).Value = 2395334.5
$ws_ALC.Cells.Item(138, 10).Value = 8549553
$ws_ALC.Cells.Item(138, 11).Value = 7186003.5
$ws_ALC.Cells.Item(138, 12).Value = 25648659
$ws_ALC.Cells.Item(138, 13).Value = -7180863.5
$ws_ALC.Cells.Item(138, 14).Value = -25658939

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Cells.Item(2, 8).Value = 47807.684
$ws_ARM.Cells.Item(2, 9).Value = 86617.086
$ws_ARM.Cells.Item(2, 10).Value = 1236.4
$ws_ARM.Cells.Item(2, 11).Value = 86617.086
$ws_ARM.Cells.Item(2, 12).Value = 1236.4
$ws_ARM.Cells.Item(2, 13).Value = -86504.086
$ws_ARM.Cells.Item(2, 14).Value = -1462.4
$ws_ARM.Cells.Item(4, 8).Value = 295.5
$ws_ARM.Cells.Item(4, 9).Value = 376.66666
$ws_ARM.Cells.Item(4, 10).Value = 52
$ws_ARM.Cells.Item(4, 11).Value = 376.66666
$ws_ARM.Cells.Item(4, 12).Value = 52
$ws_ARM.Cells.Item(4, 13).Value = -260.66666
$ws_ARM.Cells.Item(4, 14).Value = -284
$ws_ARM.Cells.Item(5, 8).Value = 62875.562
$ws_ARM.Cells.Item(5, 9).Value = 83676.586
$ws_ARM.Cells.Item(5, 10).Value = 472.5
$ws_ARM.Cells.Item(5, 11).Value = 83676.586
$ws_ARM.Cells.Item(5, 12).Value = 472.5
$ws_ARM.Cells.Item(5, 13).Value = -83564.586
$ws_ARM.Cells.Item(5, 14).Value = -696.5
$ws_ARM.Cells.Item(17, 8).Value = 0
$ws_ARM.Cells.Item(17, 10).Value = 0
$ws_ARM.Cells.Item(17, 12).Value = 0
$ws_ARM.Cells.Item(17, 14).ClearContents()
$ws_ARM.Cells.Item(41, 8).Value = 50000
$ws_ARM.Cells.Item(41, 9).Value = 0
$ws_ARM.Cells.Item(41, 10).Value = 50000
$ws_ARM.Cells.Item(41, 11).Value = 0
$ws_ARM.Cells.Item(41, 12).Value = 50000
$ws_ARM.Cells.Item(41, 13).ClearContents()
$ws_ARM.Cells.Item(41, 14).Value = -50828
$ws_ARM.Cells.Item(45, 8).Value = 1324.6364
$ws_ARM.Cells.Item(45, 9).Value = 1285.8889
$ws_ARM.Cells.Item(45, 10).Value = 1499
$ws_ARM.Cells.Item(45, 11).Value = 1285.8889
$ws_ARM.Cells.Item(45, 12).Value = 1499
$ws_ARM.Cells.Item(45, 13).Value = -908.8888999999999
$ws_ARM.Cells.Item(45, 14).Value = -2253
$ws_ARM.Cells.Item(116, 8).Value = 47807.684
$ws_ARM.Cells.Item(116, 9).Value = 86617.086
$ws_ARM.Cells.Item(116, 10).Value = 1236.4
$ws_ARM.Cells.Item(116, 11).Value = 86617.086
$ws_ARM.Cells.Item(116, 12).Value = 1236.4
$ws_ARM.Cells.Item(116, 13).Value = -84323.086
$ws_ARM.Cells.Item(116, 14).Value = -5824.4
$ws_ARM.Cells.Item(132, 8).Value = 2868.9783
$ws_ARM.Cells.Item(132, 9).Value = 2374.1843
$ws_ARM.Cells.Item(132, 11).Value = 7122.5529
$ws_ARM.Cells.Item(132, 13).Value = -4592.5529
$ws_ARM.Cells.Item(133, 8).Value = 59279.8
$ws_ARM.Cells.Item(133, 10).Value = 59279.8
$ws_ARM.Cells.Item(133, 12).Value = 59279.8
$ws_ARM.Cells.Item(133, 14).Value = -64339.8
$ws_ARM.Cells.Item(139, 8).Value = 53828.8
$ws_ARM.Cells.Item(139, 10).Value = 53828.8
$ws_ARM.Cells.Item(139, 12).Value = 53828.8
$ws_ARM.Cells.Item(139, 14).Value = -64108.8

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Cells.Item(3, 8).Value = 47807.684
$ws_BSM.Cells.Item(3, 9).Value = 86617.086
$ws_BSM.Cells.Item(3, 10).Value = 1236.4
$ws_BSM.Cells.Item(3, 11).Value = 86617.086
$ws_BSM.Cells.Item(3, 12).Value = 1236.4
$ws_BSM.Cells.Item(3, 13).Value = -86503.086
$ws_BSM.Cells.Item(3, 14).Value = -1464.4
$ws_BSM.Cells.Item(4, 8).Value = 62875.562
$ws_BSM.Cells.Item(4, 9).Value = 83676.586
$ws_BSM.Cells.Item(4, 10).Value = 472.5
$ws_BSM.Cells.Item(4, 11).Value = 83676.586
$ws_BSM.Cells.Item(4, 12).Value = 472.5
$ws_BSM.Cells.Item(4, 13).Value = -83561.586
$ws_BSM.Cells.Item(4, 14).Value = -702.5
$ws_BSM.Cells.Item(59, 8).Value = 0
$ws_BSM.Cells.Item(59, 10).Value = 0
$ws_BSM.Cells.Item(59, 12).Value = 0
$ws_BSM.Cells.Item(59, 14).ClearContents()
$ws_BSM.Cells.Item(133, 8).Value = 44000
$ws_BSM.Cells.Item(133, 10).Value = 44000
$ws_BSM.Cells.Item(133, 12).Value = 44000
$ws_BSM.Cells.Item(133, 14).Value = -54120

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Cells.Item(7, 8).Value = 73.333336
$ws_CRP.Cells.Item(7, 9).Value = 73.333336
$ws_CRP.Cells.Item(7, 10).Value = 0
$ws_CRP.Cells.Item(7, 11).Value = 73.333336
$ws_CRP.Cells.Item(7, 12).Value = 0
$ws_CRP.Cells.Item(7, 13).Value = 39.666664
$ws_CRP.Cells.Item(7, 14).ClearContents()
$ws_CRP.Cells.Item(31, 8).Value = 4503.5
$ws_CRP.Cells.Item(31, 9).Value = 0
$ws_CRP.Cells.Item(31, 10).Value = 4503.5
$ws_CRP.Cells.Item(31, 11).Value = 0
$ws_CRP.Cells.Item(31, 12).Value = 4503.5
$ws_CRP.Cells.Item(31, 13).ClearContents()
$ws_CRP.Cells.Item(31, 14).Value = -5093.5
$ws_CRP.Cells.Item(34, 8).Value = 4503.5
$ws_CRP.Cells.Item(34, 9).Value = 0
$ws_CRP.Cells.Item(34, 10).Value = 4503.5
$ws_CRP.Cells.Item(34, 11).Value = 0
$ws_CRP.Cells.Item(34, 12).Value = 4503.5
$ws_CRP.Cells.Item(34, 13).ClearContents()
$ws_CRP.Cells.Item(34, 14).Value = -4907.5
$ws_CRP.Cells.Item(58, 8).Value = 1758.9231
$ws_CRP.Cells.Item(58, 9).Value = 1091.619
$ws_CRP.Cells.Item(58, 11).Value = 1091.619
$ws_CRP.Cells.Item(58, 13).Value = -888.6189999999999
$ws_CRP.Cells.Item(88, 8).Value = 26447.5
$ws_CRP.Cells.Item(88, 10).Value = 26447.5
$ws_CRP.Cells.Item(88, 12).Value = 26447.5
$ws_CRP.Cells.Item(88, 14).Value = -27259.5
$ws_CRP.Cells.Item(91, 8).Value = 26447.5
$ws_CRP.Cells.Item(91, 10).Value = 26447.5
$ws_CRP.Cells.Item(91, 12).Value = 26447.5
$ws_CRP.Cells.Item(91, 14).Value = -29255.5
$ws_CRP.Cells.Item(105, 8).Value = 798.65216
$ws_CRP.Cells.Item(105, 9).Value = 781.1177
$ws_CRP.Cells.Item(105, 10).Value = 848.3333
$ws_CRP.Cells.Item(105, 11).Value = 781.1177
$ws_CRP.Cells.Item(105, 12).Value = 848.3333
$ws_CRP.Cells.Item(105, 13).Value = 965.8823
$ws_CRP.Cells.Item(105, 14).Value = -4342.3333
$ws_CRP.Cells.Item(107, 8).Value = 418.125
$ws_CRP.Cells.Item(107, 9).Value = 251.75
$ws_CRP.Cells.Item(107, 10).Value = 584.5
$ws_CRP.Cells.Item(107, 11).Value = 251.75
$ws_CRP.Cells.Item(107, 12).Value = 584.5
$ws_CRP.Cells.Item(107, 13).Value = 1668.25
$ws_CRP.Cells.Item(107, 14).Value = -4424.5
$ws_CRP.Cells.Item(122, 8).Value = 1709.2
$ws_CRP.Cells.Item(122, 9).Value = 1448.6666
$ws_CRP.Cells.Item(122, 11).Value = 4345.9998
$ws_CRP.Cells.Item(122, 13).Value = -1895.9998
$ws_CRP.Cells.Item(132, 8).Value = 2725.8293
$ws_CRP.Cells.Item(132, 9).Value = 2389.1035
$ws_CRP.Cells.Item(132, 10).Value = 3539.5833
$ws_CRP.Cells.Item(132, 11).Value = 7167.310500000001
$ws_CRP.Cells.Item(132, 12).Value = 10618.7499
$ws_CRP.Cells.Item(132, 13).Value = -4637.310500000001
$ws_CRP.Cells.Item(132, 14).Value = -15678.7499
$ws_CRP.Cells.Item(134, 8).Value = 3447.6316
$ws_CRP.Cells.Item(134, 9).Value = 1441.4166
$ws_CRP.Cells.Item(134, 10).Value = 6886.857
$ws_CRP.Cells.Item(134, 11).Value = 4324.2498
$ws_CRP.Cells.Item(134, 12).Value = 20660.571
$ws_CRP.Cells.Item(134, 13).Value = -1789.2498
$ws_CRP.Cells.Item(134, 14).Value = -25730.571
$ws_CRP.Cells.Item(136, 8).Value = 1758.9231
$ws_CRP.Cells.Item(136, 9).Value = 1091.619
$ws_CRP.Cells.Item(136, 11).Value = 3274.857
$ws_CRP.Cells.Item(136, 13).Value = -724.857

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Cells.Item(58, 8).Value = 7953
$ws_CUL.Cells.Item(75, 8).Value = 1996.25
$ws_CUL.Cells.Item(75, 9).Value = 1000
$ws_CUL.Cells.Item(75, 10).Value = 2138.5715
$ws_CUL.Cells.Item(75, 11).Value = 3000
$ws_CUL.Cells.Item(75, 12).Value = 6415.7145
$ws_CUL.Cells.Item(75, 13).Value = -2002
$ws_CUL.Cells.Item(75, 14).Value = -8411.7145
$ws_CUL.Cells.Item(78, 8).Value = 1996.25
$ws_CUL.Cells.Item(78, 9).Value = 1000
$ws_CUL.Cells.Item(78, 10).Value = 2138.5715
$ws_CUL.Cells.Item(78, 11).Value = 9000
$ws_CUL.Cells.Item(78, 12).Value = 19247.1435
$ws_CUL.Cells.Item(78, 13).Value = -4008
$ws_CUL.Cells.Item(78, 14).Value = -29231.1435
$ws_CUL.Cells.Item(98, 8).Value = 333.16666
$ws_CUL.Cells.Item(98, 9).Value = 333.16666
$ws_CUL.Cells.Item(98, 10).Value = 0
$ws_CUL.Cells.Item(98, 11).Value = 999.4999799999999
$ws_CUL.Cells.Item(98, 12).Value = 0
$ws_CUL.Cells.Item(98, 13).Value = 498.5000200000001
$ws_CUL.Cells.Item(98, 14).ClearContents()
$ws_CUL.Cells.Item(107, 8).Value = 465.35294
$ws_CUL.Cells.Item(107, 9).Value = 499
$ws_CUL.Cells.Item(107, 10).Value = 417.2857
$ws_CUL.Cells.Item(107, 11).Value = 1497
$ws_CUL.Cells.Item(107, 12).Value = 1251.8571
$ws_CUL.Cells.Item(107, 13).Value = 423
$ws_CUL.Cells.Item(107, 14).Value = -5091.8571
$ws_CUL.Cells.Item(122, 8).Value = 780.4375
$ws_CUL.Cells.Item(122, 9).Value = 299.83334
$ws_CUL.Cells.Item(122, 10).Value = 1068.8
$ws_CUL.Cells.Item(122, 11).Value = 2698.50006
$ws_CUL.Cells.Item(122, 12).Value = 9619.199999999999
$ws_CUL.Cells.Item(122, 13).Value = -248.5000600000003
$ws_CUL.Cells.Item(122, 14).Value = -14519.2

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Cells.Item(123, 8).Value = 10103.765
$ws_GSM.Cells.Item(123, 10).Value = 10103.765
$ws_GSM.Cells.Item(123, 12).Value = 10103.765
$ws_GSM.Cells.Item(123, 14).Value = -15003.765
$ws_GSM.Cells.Item(134, 8).Value = 7909.3335
$ws_GSM.Cells.Item(134, 10).Value = 7909.3335
$ws_GSM.Cells.Item(134, 12).Value = 23728.0005
$ws_GSM.Cells.Item(134, 14).Value = -28798.0005
$ws_GSM.Cells.Item(138, 8).Value = 68644.27
$ws_GSM.Cells.Item(138, 10).Value = 68644.27
$ws_GSM.Cells.Item(138, 12).Value = 68644.27
$ws_GSM.Cells.Item(138, 14).Value = -78924.27

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Cells.Item(19, 8).Value = 17235.334
$ws_LTW.Cells.Item(19, 10).Value = 0
$ws_LTW.Cells.Item(19, 12).Value = 0
$ws_LTW.Cells.Item(19, 14).ClearContents()
$ws_LTW.Cells.Item(46, 8).Value = 1060.125
$ws_LTW.Cells.Item(46, 10).Value = 1096.2
$ws_LTW.Cells.Item(46, 12).Value = 1096.2
$ws_LTW.Cells.Item(46, 14).Value = -1472.2
$ws_LTW.Cells.Item(55, 8).Value = 489.18518
$ws_LTW.Cells.Item(55, 9).Value = 386.25
$ws_LTW.Cells.Item(55, 10).Value = 571.5333000000001
$ws_LTW.Cells.Item(55, 11).Value = 386.25
$ws_LTW.Cells.Item(55, 12).Value = 571.5333000000001
$ws_LTW.Cells.Item(55, 13).Value = -213.25
$ws_LTW.Cells.Item(55, 14).Value = -917.5333000000001

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Cells.Item(12, 8).Value = 14624.75
$ws_WVR.Cells.Item(12, 10).Value = 2833
$ws_WVR.Cells.Item(12, 12).Value = 2833
$ws_WVR.Cells.Item(12, 14).Value = -3117
$ws_WVR.Cells.Item(24, 8).Value = 16671833
$ws_WVR.Cells.Item(24, 9).Value = 50000500
$ws_WVR.Cells.Item(24, 11).Value = 50000500
$ws_WVR.Cells.Item(24, 13).Value = -50000270
$ws_WVR.Cells.Item(82, 8).Value = 40910.57
$ws_WVR.Cells.Item(82, 10).Value = 39350.168
$ws_WVR.Cells.Item(82, 12).Value = 39350.168
$ws_WVR.Cells.Item(82, 14).Value = -40116.168
$ws_WVR.Cells.Item(85, 8).Value = 40910.57
$ws_WVR.Cells.Item(85, 10).Value = 39350.168
$ws_WVR.Cells.Item(85, 12).Value = 39350.168
$ws_WVR.Cells.Item(85, 14).Value = -42002.168
